# Applies the weekly data refresh for the Espárragos (Hortaliza) sheet.
# Each data row (2-19, skipping row 11 which is unchanged) is updated with
# new Fecha/Volumen/Precio values, reflecting the latest weekly price report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was row 16's data)
$ws.Range("D2").Value = 44875
$ws.Range("J2").Value = 300
$ws.Range("L2").Value = 1600
$ws.Range("M2").Value = 1550
$ws.Range("P2").Value = 1550

# Row 3 (was row 4's data)
$ws.Range("D3").Value = 44839
$ws.Range("J3").Value = 500
$ws.Range("K3").Value = 1700
$ws.Range("L3").Value = 1800
$ws.Range("M3").Value = 1760
$ws.Range("P3").Value = 1760

# Row 4 (was row 12's data)
$ws.Range("D4").Value = 44545
$ws.Range("J4").Value = 550
$ws.Range("M4").Value = 1755
$ws.Range("P4").Value = 1755

# Row 5 (was row 15's data)
$ws.Range("D5").Value = 44510
$ws.Range("J5").Value = 600
$ws.Range("K5").Value = 1300
$ws.Range("L5").Value = 1400
$ws.Range("M5").Value = 1350
$ws.Range("P5").Value = 1350

# Row 6 (was row 18's data)
$ws.Range("D6").Value = 44876
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 350
$ws.Range("K6").Value = 1500
$ws.Range("L6").Value = 1600
$ws.Range("M6").Value = 1557
$ws.Range("P6").Value = 1557

# Row 7 (was row 17's data)
$ws.Range("D7").Value = 44526
$ws.Range("J7").Value = 100
$ws.Range("O7").Value = "Provincia de Linares"

# Row 8 (was row 7's data)
$ws.Range("D8").Value = 44524
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("J8").Value = 200
$ws.Range("K8").Value = 1500
$ws.Range("L8").Value = 1600
$ws.Range("M8").Value = 1550
$ws.Range("O8").Value = "Provincia de Talca"
$ws.Range("P8").Value = 1550

# Row 9 (was row 14's data)
$ws.Range("D9").Value = 44519
$ws.Range("J9").Value = 250
$ws.Range("M9").Value = 1240
$ws.Range("O9").Value = "Provincia de Linares"
$ws.Range("P9").Value = 1240

# Row 10 (was row 3's data)
$ws.Range("D10").Value = 44511
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 600
$ws.Range("K10").Value = 1300
$ws.Range("L10").Value = 1400
$ws.Range("M10").Value = 1350
$ws.Range("O10").Value = "Provincia de Linares"
$ws.Range("P10").Value = 1350

# Row 12 (was row 19's data)
$ws.Range("D12").Value = 44489
$ws.Range("J12").Value = 600
$ws.Range("K12").Value = 1400
$ws.Range("L12").Value = 1500
$ws.Range("M12").Value = 1450
$ws.Range("P12").Value = 1450

# Row 13 (was row 8's data)
$ws.Range("D13").Value = 44468
$ws.Range("H13").Value = "Verde"
$ws.Range("J13").Value = 500
$ws.Range("K13").Value = 1800
$ws.Range("M13").Value = 1920
$ws.Range("N13").Value = "$/kilo"
$ws.Range("P13").Value = 1920

# Row 14 (was row 9's data)
$ws.Range("D14").Value = 44868
$ws.Range("J14").Value = 1000
$ws.Range("M14").Value = 1250
$ws.Range("O14").Value = "Región del Maule"
$ws.Range("P14").Value = 1250

# Row 15 (was row 10's data)
$ws.Range("D15").Value = 44868
$ws.Range("I15").Value = "Segunda"
$ws.Range("J15").Value = 200
$ws.Range("K15").Value = 1000
$ws.Range("L15").Value = 1000
$ws.Range("M15").Value = 1000
$ws.Range("O15").Value = "Región del Maule"
$ws.Range("P15").Value = 1000

# Row 16 (was row 5's data)
$ws.Range("D16").Value = 44881
$ws.Range("J16").Value = 200
$ws.Range("K16").Value = 2600
$ws.Range("L16").Value = 2700
$ws.Range("M16").Value = 2650
$ws.Range("P16").Value = 2650

# Row 17 (was row 6's data)
$ws.Range("D17").Value = 44881
$ws.Range("I17").Value = "Segunda"
$ws.Range("K17").Value = 2400
$ws.Range("L17").Value = 2400
$ws.Range("M17").Value = 2400
$ws.Range("P17").Value = 2400

# Row 18 (was row 2's data)
$ws.Range("D18").Value = 44860
$ws.Range("J18").Value = 1100
$ws.Range("L18").Value = 1700
$ws.Range("M18").Value = 1609
$ws.Range("P18").Value = 1609

# Row 19 (was row 13's data)
$ws.Range("D19").Value = 44496
$ws.Range("J19").Value = 550
$ws.Range("K19").Value = 1500
$ws.Range("L19").Value = 2000
$ws.Range("M19").Value = 1773
$ws.Range("N19").Value = "$/paquete"
$ws.Range("P19").Value = 1773

